# Update countries & provincias Spain
# - refresh the "last updated" timestamp
# - refresh case counters for several countries (new data pulled in)
# - re-sort a few rows whose total-case count changed enough to change
#   their rank in the (descending, by "Casos totales") table:
#     Cuba/Oman/Honduras/Uzbekistan block and the Nigeria/Malta/Ghana block
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 5 de Abril de 2020 a las 08:22'
$ws.Cells.Item(4, 2).Value = 311637
$ws.Cells.Item(4, 3).Value = 280
$ws.Cells.Item(4, 5).Value = 288358
$ws.Cells.Item(21, 4).Value = 477
$ws.Cells.Item(21, 5).Value = 7495
$ws.Cells.Item(66, 2).Value = 811
$ws.Cells.Item(66, 3).Value = 40
$ws.Cells.Item(66, 5).Value = 793
$ws.Cells.Item(75, 5).Value = 509
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = 6
$ws.Cells.Item(88, 2).Value = 363
$ws.Cells.Item(88, 3).Value = 8
$ws.Cells.Item(88, 4).Value = 54
$ws.Cells.Item(88, 5).Value = 304
$ws.Cells.Item(89, 4).Value = 15
$ws.Cells.Item(89, 5).Value = 315
$ws.Cells.Item(94, 1).Value = 'Uzbekistan'
$ws.Cells.Item(94, 2).Value = 298
$ws.Cells.Item(94, 3).Value = 32
$ws.Cells.Item(94, 4).Value = 25
$ws.Cells.Item(94, 5).Value = 271
$ws.Cells.Item(94, 6).Value = 8
$ws.Cells.Item(94, 8).Value = 2
$ws.Cells.Item(95, 2).Value = 298
$ws.Cells.Item(95, 3).Value = 21
$ws.Cells.Item(95, 5).Value = 235
$ws.Cells.Item(96, 1).Value = 'Cuba'
$ws.Cells.Item(96, 2).Value = 288
$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 4).Value = 15
$ws.Cells.Item(96, 5).Value = 267
$ws.Cells.Item(96, 6).Value = 11
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 6
$ws.Cells.Item(97, 1).Value = 'Honduras'
$ws.Cells.Item(97, 2).Value = 268
$ws.Cells.Item(97, 3).Value = 4
$ws.Cells.Item(97, 4).Value = 6
$ws.Cells.Item(97, 5).Value = 240
$ws.Cells.Item(97, 6).Value = 10
$ws.Cells.Item(97, 7).Value = 7
$ws.Cells.Item(97, 8).Value = 22
$ws.Cells.Item(104, 1).Value = 'Ghana'
$ws.Cells.Item(104, 2).Value = 214
$ws.Cells.Item(104, 3).Value = 9
$ws.Cells.Item(104, 4).Value = 31
$ws.Cells.Item(104, 5).Value = 178
$ws.Cells.Item(104, 6).Value = 2
$ws.Cells.Item(104, 8).Value = 5
$ws.Cells.Item(105, 1).Value = 'Malta'
$ws.Cells.Item(105, 2).Value = 213
$ws.Cells.Item(105, 4).Value = 2
$ws.Cells.Item(105, 5).Value = 211
$ws.Cells.Item(105, 6).Value = 3
$ws.Cells.Item(105, 8).Value = 0
